$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 33
$ws.Range("H33").Value = 92820.164
$ws.Range("I33").Value = 139030.62
$ws.Range("J33").Value = 399.25
$ws.Range("K33").Value = 139030.62
$ws.Range("L33").Value = 399.25
$ws.Range("M33").Value = -138801.62
$ws.Range("N33").Value = -857.25

# row 62
$ws.Range("H62").Value = 2839.889
$ws.Range("I62").Value = 2362.3333
$ws.Range("J62").Value = 3795
$ws.Range("K62").Value = 2362.3333
$ws.Range("L62").Value = 3795
$ws.Range("M62").Value = -1738.3333
$ws.Range("N62").Value = -5043

# row 65
$ws.Range("H65").Value = 2839.889
$ws.Range("I65").Value = 2362.3333
$ws.Range("J65").Value = 3795
$ws.Range("K65").Value = 11811.6665
$ws.Range("L65").Value = 18975
$ws.Range("M65").Value = -8691.666499999999
$ws.Range("N65").Value = -25215

# row 74
$ws.Range("H74").Value = 3001496.8
$ws.Range("I74").Value = 3001496.8
$ws.Range("K74").Value = 3001496.8
$ws.Range("M74").Value = -3000560.8

# row 77
$ws.Range("H77").Value = 3001496.8
$ws.Range("I77").Value = 3001496.8
$ws.Range("K77").Value = 15007484
$ws.Range("M77").Value = -15002804

# row 103
$ws.Range("H103").Value = 1743.375
$ws.Range("J103").Value = 1225
$ws.Range("L103").Value = 3675
$ws.Range("N103").Value = -4847

# row 137
$ws.Range("H137").Value = 2579.2856
$ws.Range("I137").Value = 2175.889
$ws.Range("K137").Value = 6527.667
$ws.Range("M137").Value = -3977.667

$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value = 1893.6938
$ws.Range("I32").Value = 1861.8445
$ws.Range("K32").Value = 1861.8445
$ws.Range("M32").Value = -1574.8445

# row 46
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").ClearContents()
$ws.Range("N46").Value = 0

# row 76
$ws.Range("H76").Value = 30997
$ws.Range("J76").Value = 30997
$ws.Range("L76").Value = 30997
$ws.Range("N76").Value = -31673

# row 79
$ws.Range("H79").Value = 30997
$ws.Range("J79").Value = 30997
$ws.Range("L79").Value = 30997
$ws.Range("N79").Value = -33337

# row 110
$ws.Range("H110").Value = 201349.6
$ws.Range("I110").Value = 201349.6
$ws.Range("K110").Value = 201349.6
$ws.Range("M110").Value = -199304.6

$ws = $wb.Worksheets.Item("BSM")
# row 86
$ws.Range("H86").Value = 3342.6667
$ws.Range("J86").Value = 3501.9092
$ws.Range("L86").Value = 3501.9092
$ws.Range("N86").Value = -5747.9092

# row 89
$ws.Range("H89").Value = 3342.6667
$ws.Range("J89").Value = 3501.9092
$ws.Range("L89").Value = 17509.546
$ws.Range("N89").Value = -28741.546

# row 105
$ws.Range("H105").Value = 1935
$ws.Range("I105").Value = 2003.1111
$ws.Range("J105").Value = 1781.75
$ws.Range("K105").Value = 2003.1111
$ws.Range("L105").Value = 1781.75
$ws.Range("M105").Value = -256.1111000000001
$ws.Range("N105").Value = -5275.75

$ws = $wb.Worksheets.Item("CRP")
# row 58
$ws.Range("H58").Value = 26325032
$ws.Range("I58").Value = 31260426
$ws.Range("J58").Value = 2932.6667
$ws.Range("K58").Value = 31260426
$ws.Range("L58").Value = 2932.6667
$ws.Range("M58").Value = -31260223
$ws.Range("N58").Value = -3338.6667

# row 62
$ws.Range("H62").Value = 3231.4285
$ws.Range("J62").Value = 3231.4285
$ws.Range("L62").Value = 3231.4285
$ws.Range("N62").Value = -4479.4285

# row 65
$ws.Range("H65").Value = 3231.4285
$ws.Range("J65").Value = 3231.4285
$ws.Range("L65").Value = 16157.1425
$ws.Range("N65").Value = -22397.1425

# row 88
$ws.Range("H88").Value = 18500
$ws.Range("I88").Value = 15000
$ws.Range("J88").Value = 20600
$ws.Range("K88").Value = 15000
$ws.Range("L88").Value = 20600
$ws.Range("M88").Value = -14594
$ws.Range("N88").Value = -21412

# row 91
$ws.Range("H91").Value = 18500
$ws.Range("I91").Value = 15000
$ws.Range("J91").Value = 20600
$ws.Range("K91").Value = 15000
$ws.Range("L91").Value = 20600
$ws.Range("M91").Value = -13596
$ws.Range("N91").Value = -23408

# row 99
$ws.Range("H99").Value = 3000
$ws.Range("I99").Value = 3000
$ws.Range("K99").Value = 3000
$ws.Range("M99").Value = -1502

# row 112
$ws.Range("H112").Value = 90000
$ws.Range("J112").Value = 90000
$ws.Range("L112").Value = 90000
$ws.Range("N112").Value = -92954

# row 122
$ws.Range("H122").Value = 2374.375
$ws.Range("I122").Value = 2374.375
$ws.Range("K122").Value = 7123.125
$ws.Range("M122").Value = -4673.125

# row 126
$ws.Range("H126").Value = 3000
$ws.Range("I126").Value = 3000
$ws.Range("K126").Value = 9000
$ws.Range("M126").Value = -6530

# row 132
$ws.Range("H132").Value = 100003840
$ws.Range("I132").Value = 125003800
$ws.Range("J132").Value = 4007
$ws.Range("K132").Value = 375011400
$ws.Range("L132").Value = 12021
$ws.Range("M132").Value = -375008870
$ws.Range("N132").Value = -17081

# row 136
$ws.Range("H136").Value = 26325032
$ws.Range("I136").Value = 31260426
$ws.Range("J136").Value = 2932.6667
$ws.Range("K136").Value = 93781278
$ws.Range("L136").Value = 8798.000100000001
$ws.Range("M136").Value = -93778728
$ws.Range("N136").Value = -13898.0001

$ws = $wb.Worksheets.Item("CUL")
# row 4
$ws.Range("H4").Value = 7251760
$ws.Range("I4").Value = 8924089
$ws.Range("K4").Value = 26772267
$ws.Range("M4").Value = -26772155

# row 92
$ws.Range("H92").Value = 538.6667
$ws.Range("J92").Value = 546.4
$ws.Range("L92").Value = 1639.2
$ws.Range("N92").Value = -4135.2

# row 131
$ws.Range("H131").Value = 1989.625
$ws.Range("J131").Value = 2247
$ws.Range("L131").Value = 6741
$ws.Range("N131").Value = -16821

$ws = $wb.Worksheets.Item("GSM")
# row 80
$ws.Range("H80").Value = 2552.7144
$ws.Range("I80").Value = 2650.1667
$ws.Range("J80").Value = 2479.625
$ws.Range("K80").Value = 2650.1667
$ws.Range("L80").Value = 2479.625
$ws.Range("M80").Value = -1652.1667
$ws.Range("N80").Value = -4475.625

# row 83
$ws.Range("H83").Value = 2552.7144
$ws.Range("I83").Value = 2650.1667
$ws.Range("J83").Value = 2479.625
$ws.Range("K83").Value = 13250.8335
$ws.Range("L83").Value = 12398.125
$ws.Range("M83").Value = -8258.833500000001
$ws.Range("N83").Value = -22382.125

# row 122
$ws.Range("H122").Value = 4601.5884
$ws.Range("I122").Value = 3902.6365
$ws.Range("J122").Value = 5883
$ws.Range("K122").Value = 11707.9095
$ws.Range("L122").Value = 17649
$ws.Range("M122").Value = -9257.9095
$ws.Range("N122").Value = -22549

# row 126
$ws.Range("H126").Value = 9398.799999999999
$ws.Range("I126").Value = 3000
$ws.Range("J126").Value = 18997
$ws.Range("K126").Value = 9000
$ws.Range("L126").Value = 56991
$ws.Range("M126").Value = -6530
$ws.Range("N126").Value = -61931

$ws = $wb.Worksheets.Item("LTW")
# row 34
$ws.Range("H34").Value = 5000
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()

# row 54
$ws.Range("H54").Value = 30000
$ws.Range("I54").Value = 30000
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 30000
$ws.Range("L54").ClearContents()
$ws.Range("N54").Value = 0
$ws.Range("M54").Value = -29356

# row 68
$ws.Range("H68").Value = 1403.8
$ws.Range("I68").Value = 1198.2858
$ws.Range("J68").Value = 1883.3334
$ws.Range("K68").Value = 1198.2858
$ws.Range("L68").Value = 1883.3334
$ws.Range("M68").Value = -449.2858000000001
$ws.Range("N68").Value = -3381.3334

# row 71
$ws.Range("H71").Value = 1403.8
$ws.Range("I71").Value = 1198.2858
$ws.Range("J71").Value = 1883.3334
$ws.Range("K71").Value = 5991.429
$ws.Range("L71").Value = 9416.666999999999
$ws.Range("M71").Value = -2247.429
$ws.Range("N71").Value = -16904.667

# row 136
$ws.Range("H136").Value = 1788.0526
$ws.Range("I136").Value = 1508.7273
$ws.Range("J136").Value = 2172.125
$ws.Range("K136").Value = 4526.1819
$ws.Range("L136").Value = 6516.375
$ws.Range("M136").Value = -1976.1819
$ws.Range("N136").Value = -11616.375

$ws = $wb.Worksheets.Item("WVR")
# row 82
$ws.Range("H82").Value = 20421.6
$ws.Range("J82").Value = 20421.6
$ws.Range("L82").Value = 20421.6
$ws.Range("N82").Value = -21187.6

# row 85
$ws.Range("H85").Value = 20421.6
$ws.Range("J85").Value = 20421.6
$ws.Range("L85").Value = 20421.6
$ws.Range("N85").Value = -23073.6

# row 132
$ws.Range("H132").Value = 20004300
$ws.Range("I132").Value = 25002300
$ws.Range("K132").Value = 75006900
$ws.Range("M132").Value = -75004370
